# Apply stricter validation formatting fixes to the "variable_mapping" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variable_mapping")

# 1. Rename the two header-label strings (remove the embedded spaces).
$ws.Range("A1").Value = "PlatformName"
$ws.Range("B1").Value = "HeaderRow"

# 2. Unify the style of columns C:O with column A: left-aligned text instead
#    of the previously separate "general" / duplicate "left" styles.
$ws.Range("C1:O7").HorizontalAlignment = -4131  # xlLeft
$ws.Columns.Item(1).HorizontalAlignment = -4131 # xlLeft (column A, general -> left)

# 3. Bump the row height of the data rows (2-7) slightly, leaving the header
#    row (1) untouched.
for ($r = 2; $r -le 7; $r++) {
    $ws.Rows.Item($r).RowHeight = 20.25
}
